$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 - new time log entry (Architecture / Continued work on UML class models.)
$ws.Range("A6").Value = 44495
$ws.Range("B6").Value = 20/24
$ws.Range("C6").Value = 21/24
$ws.Range("E6").Value = "Architecture"
$ws.Range("F6").Value = "Continued work on UML class models."

# Row 7 - new time log entry (Code / Writing code from the UML diagrams.)
$ws.Range("B7").Value = 21/24
$ws.Range("C7").Value = (21*60+57)/1440
$ws.Range("E7").Value = "Code"
$ws.Range("F7").Value = "Writing code from the UML diagrams."
